$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Expand the first bullet under "Network Analytics Engine" from the
#    placeholder "Point1" into the full description, moving the
#    document's "_GoBack" (last-edit) bookmark into the new text at the
#    point where the author paused while typing.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Point1")
$point1Start = $rng.Start
$fullText = "Software Engineer for a team working on Network Analytics Engine module on a network device such as switch."
$rng.Text = $fullText

# The author's cursor (and so the "_GoBack" bookmark Word drops at the
# last edited spot) ended up right after "...Engine mo", i.e. just
# before "dule on a network...".
$beforeBookmark = "Software Engineer for a team working on Network Analytics Engine mo"
$bookmarkPos = $point1Start + $beforeBookmark.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# ---------------------------------------------------------------------
# 2. Collapse the " " + "Python" runs (Environment: <bold> Python) into
#    a single run.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Python", $true, $false, $false, $false, $false,
    $true, 1, $false, "Python", 2)

# ---------------------------------------------------------------------
# 3. Merge "Graduate Assistant" + "," into one run. A temporary bookmark
#    is used to stop the surrounding same-formatted runs from also
#    collapsing into this one; it is removed immediately afterwards.
# ---------------------------------------------------------------------
$gaRng = $d.Content
$gaRng.Find.Execute("Graduate Assistant,")
$gaEnd = $gaRng.End
$gaBoundary = $d.Range($gaEnd, $gaEnd)
$d.Bookmarks.Add("tempBoundaryGA", $gaBoundary)

$d.Content.Find.Execute(
    "Graduate Assistant,", $true, $false, $false, $false, $false,
    $true, 1, $false, "Graduate Assistant,", 2)

$d.Bookmarks("tempBoundaryGA").Delete()

# ---------------------------------------------------------------------
# 4. Merge "Uni" + "versity of Arizona" (which also removes the old
#    "_GoBack" bookmark that used to sit between them - it already
#    moved to step 1) into a single run reading "University of Arizona".
#    A temporary leading boundary keeps the preceding "CE, " run intact.
# ---------------------------------------------------------------------
$uaRng = $d.Content
$uaRng.Find.Execute("University of Arizona")
$uaStart = $uaRng.Start
$uaBoundary = $d.Range($uaStart, $uaStart)
$d.Bookmarks.Add("tempBoundaryUA", $uaBoundary)

$d.Content.Find.Execute(
    "University of Arizona", $true, $false, $false, $false, $false,
    $true, 1, $false, "University of Arizona", 2)

$d.Bookmarks("tempBoundaryUA").Delete()
